# "updated test data excel"
#
# The date pieces that used to read "Jan" / "31" move forward a month:
#   - OrderConfirmationPage / StoreHeadConfirmationPage: the from/to date
#     pair goes from "Jan 31" to "Feb 1".
#   - CreateCustomerOrderNo: the occasion / due-date / cheque-date triples
#     go from "Jan 31" to "Feb 28".
# A few sheets were also left scrolled/selected differently.

$wb = $excel.ActiveWorkbook

# --- OrderHeaderPage: selection only moved from H2 to J2 ---
$ws = $wb.Worksheets.Item("OrderHeaderPage")
$ws.Range("J2").Select() | Out-Null

# --- OrderDetailPage: viewport + selection moved ---
$ws = $wb.Worksheets.Item("OrderDetailPage")
$ws.Range("AM2").Select() | Out-Null

# --- OrderConfirmationPage: from/to date changed Jan 31 -> Feb 1 ---
$ws = $wb.Worksheets.Item("OrderConfirmationPage")
$ws.Range("C2").Value = "Feb"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "Feb"
$ws.Range("H2").Value = "1"
$ws.Range("I2").Select() | Out-Null

# --- StoreHeadConfirmationPage: same from/to date change (selection here
#     was already I2 both before and after) ---
$ws = $wb.Worksheets.Item("StoreHeadConfirmationPage")
$ws.Range("C2").Value = "Feb"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "Feb"
$ws.Range("H2").Value = "1"

# --- CreateCustomerOrderNo: occasion / due date / cheque date triples
#     changed Jan 31 -> Feb 28. This sheet is the active tab both before
#     and after the edit, so re-activate it (and its original A2
#     selection) last, after the other sheets' selections were touched. ---
$ws = $wb.Worksheets.Item("CreateCustomerOrderNo")
$ws.Range("H2").Value = "Feb"
$ws.Range("J2").Value = "28"
$ws.Range("W2").Value = "Feb"
$ws.Range("Y2").Value = "28"
$ws.Range("AU2").Value = "Feb"
$ws.Range("AW2").Value = "28"
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
